$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'59.358.38"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.17%  '

$ws.Range('D3').Value = "'3.004.66"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.99%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = "'564.29"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.01%  '

$ws.Range('D6').Value = "'139.38"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.24%  '

$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').Value = "'0.520"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.12%  '

$ws.Range('D9').Value = "'2.996.83"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.82%  '

$ws.Range('D10').Value = "'0.133"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.14%  '

$ws.Range('D11').Value = "'5.21"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.82%  '

$ws.Range('E12').Value = '  +2.14%  '

$ws.Range('E13').Value = '  +3.38%  '

$ws.Range('D14').Value = "'33.87"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.49%  '

$ws.Range('E15').Value = '  +2.18%  '

$ws.Range('E16').Value = '  +7.09%  '

$ws.Range('D17').Value = "'3.503.88"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.15%  '

$ws.Range('D18').Value = "'3.008.21"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.33%  '

$ws.Range('D19').Value = "'59.356.41"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.36%  '

$ws.Range('D20').Value = "'432.28"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.81%  '

$ws.Range('D21').Value = "'13.68"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.26%  '

$ws.Range('D22').Value = "'0.724"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.29%  '

$ws.Range('D23').Value = "'13.61"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.66%  '

$ws.Range('D24').Value = "'7.16"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.08%  '

$ws.Range('D25').Value = "'80.66"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.27%  '

$ws.Range('E26').Value = '  -0.03%  '

$ws.Range('E27').Value = '  +11.88%  '

$ws.Range('E28').Value = '  +0.20%  '

$ws.Range('E29').Value = '  +2.29%  '

$ws.Range('D30').Value = "'7.91"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.10%  '

$ws.Range('D31').Value = "'25.82"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.95%  '

$ws.Range('D32').Value = "'6.14"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.81%  '

$ws.Range('E33').Value = '  +0.28%  '

$ws.Range('D34').Value = "'1.00"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.68%  '

$ws.Range('E35').Value = '  +6.07%  '

$ws.Range('D36').Value = "'0.0₃0760"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.63%  '

$ws.Range('D37').Value = "'2.12"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.21%  '

$ws.Range('D38').Value = "'49.08"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.66%  '

$ws.Range('D39').Value = "'8.67"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.12%  '

$ws.Range('E40').Value = '  +6.73%  '

$ws.Range('D41').Value = "'411.43"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.07%  '

$ws.Range('D42').Value = "'0.0355"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.97%  '

$ws.Range('D43').Value = "'2.776.58"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.16%  '

$ws.Range('E44').Value = '  -0.28%  '

$ws.Range('E45').Value = '  +4.31%  '

$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = "'0.999"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.01%  '

$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = "'35.10"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +21.28%  '

$ws.Range('D48').Value = "'123.66"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.31%  '

$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = "'0.111"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.45%  '

$ws.Range('B50').Value = 'Fetch.AI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D50').Value = "'2.02"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.03%  '

$ws.Range('D51').Value = "'23.63"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.03%  '
